$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Unit Test Plan Preparation
$ws.Range("L14").Value = 1
$ws.Range("O14").Value = 1

# Row 19: LLD Rework
$ws.Range("O19").Value = 1.5

# Rows 24-31: add hours in column Q
$ws.Range("Q24").Value = 1
$ws.Range("Q25").Value = 1
$ws.Range("Q26").Value = 1
$ws.Range("Q27").Value = 2
$ws.Range("Q28").Value = 2
$ws.Range("Q29").Value = 1
$ws.Range("Q30").Value = 1
$ws.Range("Q31").Value = 1

# Rows 32-33: add hours in column R
$ws.Range("R32").Value = 3
$ws.Range("R33").Value = 1

# Update the selected cell to match the saved view state
$ws.Range("T16").Select()
